# refactor m:n field name
# Renames the "modules/user_*" / "moduler/user_*" identifiers in the
# ErrorMessages sheet so the module segment uses a dot instead of an
# underscore before the sub-module name (e.g. user_password -> user.password).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("ErrorMessages")

$ws.Range("C8").Value  = "modules/user.password.logUp"
$ws.Range("C9").Value  = "modules/user.password.logUp"

# C8:C9 previously had no explicit formatting; align them with the matching
# style used by the neighbouring cells (C10, C11, ...) in this column.
$ws.Range("C8:C9").Font.Family = 3

$ws.Range("C10").Value = "modules/user.password.logIn"
$ws.Range("C11").Value = "modules/user.password.logIn"
$ws.Range("C12").Value = "modules/user.password.logIn"
$ws.Range("C13").Value = "modules/user.password.logIn"

$ws.Range("C14").Value = "modules/user.search.all"
$ws.Range("C15").Value = "modules/user.search.choice"
$ws.Range("C16").Value = "modules/user.search.search_init"
$ws.Range("C17").Value = "modules/user.search.search_include"

$ws.Range("C18").Value = "moduler/user.update.patch_user"
$ws.Range("C19").Value = "moduler/user.update.patch_info"

# The ErrorMessages sheet becomes the active tab/selection (was API sheet).
$ws.Activate()
$ws.Range("E27").Select()
